# Update the "想去人数" (interested-people count) figures in the F column
# on both the "展览" and "全部类型" worksheets, mirroring the values that
# were regenerated for the gh-pages data refresh at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Map of F-cell row -> new value, identical for both sheets that carry the
# full exhibition table.
$updates = @{
    3  = 101
    4  = 1572
    6  = 1090
    7  = 11318
    13 = 782
    14 = 12309
    15 = 12955
    22 = 86
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
